$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("B3:B60")
$fcs = $rng.FormatConditions
$fc1 = $fcs.Item(1)
$fc1.Delete()
$t = [Microsoft.Office.Interop.Excel.XlFormatConditionType]::xlNoBlanksCondition
$newFc = $fcs.Add($t, [System.Type]::Missing, [System.Type]::Missing)
Write-Host "added" $t
